$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the comment text for the Hallmark Movies Now / YouTube TV row:
# "Network Added to Add-On Service" -> "Network Added to Add-On Package"
$ws.Range("E7").Value = "Network Added to Add-On Package"

# Move active selection to A2 (top-left cell below the frozen header row)
$ws.Range("A2").Select()
